# "still studying the parameters on excell"
# Insert a new "#splitLength" column before column J (the existing J..O data
# shifts right to K..P), fill in the new column's header + values, then
# leave the selection where the author left it (H19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at J - everything from J onward (incl. the #minPredictionRatio
# column and the rest) shifts one column to the right automatically.
$ws.Columns("J:J").Insert()

# Header for the freshly inserted column.
$ws.Range("J1").Value = "#splitLength"

# Body values for the new "#splitLength" column (rows 2-24).
$splitLengths = @(10, 20, 6, 6, 8, 22, 12, 22, 22, 22, 22, 22, 22, 22, 22, 22, 22, 22, 22, 22, 22, 22, 22)
for ($i = 0; $i -lt $splitLengths.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 10).Value = $splitLengths[$i]
}

# Match the author's last selection/cursor position.
$ws.Range("H19").Select()
